$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Name" header in A1
$ws.Range("A1").ClearContents()

# Update first-test-case first name value
$ws.Range("B2").Value = "Marco1"

# Fix typo/trailing-space in email header
$ws.Range("C1").Value = "email"

# Move active selection to D1
$ws.Range("D1").Select()
